# edit.ps1 - apply the "Robotic Assistance" -> "Biology / Cells" essay rewrite
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title / byline / email -------------------------------------------------
Replace-Text "Robotic Assistance in Medical Surgeries" "The Intricate Workings of Cells: A Journey into Biology"
Replace-Text " Sarah Harper" " Isabelle Leclerc"
Replace-Text "harper" "isabelle"
Replace-Text "sarah@deakin" "leclerc@schooldistrict"

# --- Intro paragraph sentences ----------------------------------------------
Replace-Text "The fusion of robotics with the intricacies of surgery has ignited a revolution in the healthcare landscape, promising unparalleled precision, reduced invasiveness, and enhanced patient outcomes" "Biology, the study of life, stands as a captivating field"
Replace-Text "From minimally invasive procedures to complex operations, surgical robots are transforming the way medical interventions are performed, redefining the boundaries of what was once deemed impossible" "It opens doors to understanding the intricate workings that govern living organisms, from the tiniest microscopic cells to the grandest of creatures"
Replace-Text "As the vanguard of medical technology, surgical robots are equipped with exceptional dexterity, steady precision, and unwavering accuracy, offering an unprecedented level of control to surgeons" "As we embark on this journey into the realm of biology, we will unravel the mysteries of cells, the fundamental building blocks of life"
Replace-Text "This essay delves into the compelling world of robotic assistance in medical surgeries, exploring its multifaceted applications, the benefits it bestows upon patients and healthcare professionals, and the promising future it holds" "We will explore their diverse structures, functions, and mechanisms, shedding light on the delicate harmony that sustains living systems"

# --- First body paragraph sentences -----------------------------------------
Replace-Text "Surgical robots excel in performing intricate procedures with exceptional precision, minimizing the invasiveness of surgical interventions and reducing the associated risks for patients" "Cells, the smallest units capable of independent life, present a fascinating microcosm of complexity"
Replace-Text "The robotic arms, equipped with miniaturized instruments, can navigate through narrow and delicate anatomical structures with unparalleled accuracy, enabling surgeons to access areas that would otherwise be inaccessible" "Within their minute boundaries, a symphony of biological processes takes place, governed by the delicate interplay of molecules and organelles"
Replace-Text "Moreover, the steady and precise movements of the robotic system eliminate tremors or fatigue-induced errors, ensuring a consistent level of surgical precision throughout the procedure" "These tiny powerhouses perform essential functions necessary for survival, including energy production, waste removal, and the synthesis of vital molecules"
Replace-Text "This heightened precision translates into reduced trauma to surrounding tissues, less blood loss, and a swifter recovery for patients" "By delving into the intricate world of cells, we discover the secrets of life's fundamental processes"

# --- Second body paragraph sentences -----------------------------------------
Replace-Text "The application of robotics in surgeries extends beyond precision to encompass a broad range of surgical specialties" "Furthermore, the study of cells provides a foundation for understanding the complexities of living organisms"
Replace-Text "From delicate cardiovascular procedures to complex neurosurgeries, robotic assistance has revolutionized the way surgeons approach intricate interventions" "It allows us to comprehend how cells cooperate to form tissues, which in turn combine to create organs and systems within multicellular organisms"
Replace-Text "In cardiac surgery, robotic systems facilitate minimally invasive coronary artery bypass grafting and valve repair, minimizing the invasiveness of these procedures and reducing the risks associated with open heart surgery" "This understanding underlies our comprehension of how these organisms grow, reproduce, and interact with their environment"
Replace-Text "Neurosurgeons harness the capabilities of robotic systems to remove brain tumors with greater accuracy, preserving healthy tissue and minimizing neurological complications" "By exploring the marvels of cells, we gain insights into the miraculous processes that sustain life"

# Remove the final sentence of the second body paragraph ("General surgeons ...")
# and grow the paragraph with the new "Body:" section instead.
$d.Content.Find.Execute(" General surgeons employ robotic assistance for precise gall bladder removal, hernia repair, and colorectal surgery, yielding improved patient outcomes and reduced recovery times.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$nl = [char]11
$newBody = $nl + $nl + "Body:" + $nl + $nl + `
    "Cells display a remarkable diversity, ranging from simple prokaryotes to complex eukaryotes." + `
    " Prokaryotes, such as bacteria and archaea, lack a nucleus and other membrane-bound organelles, while eukaryotes, including plants and animals, possess these advanced cellular structures." + `
    " Specialized cells, each with unique functions, contribute to the intricate organization of multicellular organisms." + `
    " Epithelial cells protect the body's surfaces, muscle cells contract to enable movement, and nerve cells transmit electrical impulses." + `
    " This cellular diversity reflects the extraordinary complexity and adaptability of life." + $nl + $nl + `
    "The processes occurring within cells are equally diverse." + `
    " Cells engage in a continuous exchange of materials with their surroundings, taking in nutrients and oxygen while releasing waste products." + `
    " They convert these nutrients into energy through respiration, a process that releases energy stored in food molecules." + `
    " Cells also synthesize proteins, the building blocks of life, using genetic information stored in DNA." + `
    " These intricate processes are essential for growth, repair, and reproduction, highlighting the dynamic nature of cellular life." + $nl + $nl + `
    "Cells communicate with each other through a variety of mechanisms, including chemical signals, electrical impulses, and direct physical contact." + `
    " This intercellular communication is crucial for coordinating cellular activities and maintaining the overall integrity of the organism." + `
    " It enables cells to respond to changes in their environment, adapt to various conditions, and collectively carry out complex functions." + `
    " The ability of cells to communicate underscores the interconnectedness and cooperative nature of life"

$para2 = $d.Paragraphs(5).Range
$insPoint = $d.Range($para2.End - 1, $para2.End - 1)
$insPoint.Text = $newBody

# --- Remove the lastRenderedPageBreak hint before "Summary" -----------------
$summaryPar = $d.Paragraphs(6).Range
$summaryPar.Text = "Summary"

# --- Summary paragraph sentences --------------------------------------------
Replace-Text "Robotic assistance in medical surgeries has heralded a new era in healthcare, transforming the way surgical interventions are performed" "Biology, the study of life, offers a profound understanding of the intricate workings of cells, the fundamental building blocks of living organisms"
Replace-Text "With their exceptional precision, reduced invasiveness, and adaptability across various surgical specialties, robotic systems have become invaluable allies to surgeons, enabling them to perform complex procedures with greater accuracy and efficiency" "From the diversity of cells to the processes occurring within them, the field of biology unveils the secrets of life's fundamental processes"
Replace-Text "The benefits of robotic surgery extend to patients, offering reduced trauma, faster recovery, and improved overall outcomes" "Exploring the world of cells provides insights into the miraculous processes that sustain life and reveals the extraordinary complexity and adaptability of living organisms"

# Remove the final two sentences and replace with one closing sentence.
$d.Content.Find.Execute(" As technology continues to advance, the integration of robotics in surgeries will only deepen, leading to even more remarkable advancements in healthcare. The future of robotic-assisted surgery is limitless, promising to revolutionize the surgical landscape further, improving patient care, and unlocking new possibilities in the pursuit of surgical excellence.", $true, $false, $false, $false, $false, $true, 1, $false, " Through the study of biology, we gain a deeper appreciation for the unity and diversity of life, fostering a greater understanding of our place within the natural world.", 2) | Out-Null

# --- Trailing empty paragraph -------------------------------------------------
$endPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endPoint.Text = [char]13

Write-Output "edit complete"
